$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-12 (column A label, B/C/D numeric values)
$data = @(
    @("Diseases (patient-stated)", 4.3, 7.6, 5.2),
    @("Injuries & adverse effects", 7.7, 12.1, 6.1),
    @("Other", 6.1, 7.3, 5.6),
    @("Symptom – Circulatory", 7.5, 12.1, 9.4),
    @("Symptom – Digestive", 9.9, 12.8, 11.6),
    @("Symptom – General", 4.6, 4.8, 5.5),
    @("Symptom – Genitourinary", 5, 6.2, 5.5),
    @("Symptom – Nervous", 11.8, 9.8, 13.3),
    @("Symptom – Respiratory", 38.9, 21.4, 35),
    @("Symptom – Skin/Hair/Nails", 2.2, 2.9, 1.9),
    @("Uncodable/Unknown", 1.9, 3.1, 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
